$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that needs to be bumped
# from 45180 (2023-09-11) to 45181 (2023-09-12) for every data row (2-170).
for ($r = 2; $r -le 170; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value = 45181
    }
}
